$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Urls"
$ws.Range("A1").Errors.Item(3).Ignore = $true
